$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Yasser009@gmail.com"
$ws.Range("D3").Value = "Yasser010@gmail.com"
$ws.Range("D4").Value = "Yasser011s@gmail.com"
$ws.Range("D5").Value = "Yasser008@gmail.com"

$ws.Range("D4").Select()
